$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(4)

# Crop the top of the picture (removes the top ~34.694% of the source image),
# matching the new <a:srcRect t="34694"/> introduced by the edit.
$sh.PictureFormat.CropTop = 37.98993

# Reposition / resize the picture frame to its new place on the slide.
$sh.Left   = 492.1526336669922
$sh.Top    = 330.25146484375
$sh.Width  = 222.39791870117188
$sh.Height = 83.15578842163086
